# Prepare to merge into main
# - Populate the "Thermo" sheet (sheet2) with its raw-score -> percentile
#   lookup table (matching the style already used on the "1110" and "1220"
#   sheets), and update the sheets' view/selection state so that "Thermo"
#   ends up the active tab.

$wb = $excel.ActiveWorkbook

$ws1110   = $wb.Worksheets.Item("1110")
$wsThermo = $wb.Worksheets.Item("Thermo")
$ws1220   = $wb.Worksheets.Item("1220")

# ---------------------------------------------------------------------
# Thermo lookup table
# ---------------------------------------------------------------------
$wsThermo.Cells.Item(1, 1).Value = "Lookup Table"
$wsThermo.Cells.Item(2, 1).Value = "Raw Score"
$wsThermo.Cells.Item(2, 2).Value = "%-ile"

$percentiles = @(100,100,100,100,100,100,99,98,97,96,95,94,92,90,88,85,83,80,76,71,66,60,55,48,41,36,31,26,21,16,12,9,7,5,4,3,2,1,0,0,0,0,0,0,0,0,0,0,0,0)

$row = 3
for ($i = 0; $i -lt $percentiles.Length; $i++) {
    $wsThermo.Cells.Item($row, 1).Value = 53 - $row
    $wsThermo.Cells.Item($row, 2).Value = $percentiles[$i]
    $row++
}

$wsThermo.Columns.Item(1).AutoFit() | Out-Null

# ---------------------------------------------------------------------
# View state / selections
# ---------------------------------------------------------------------
# "1110" keeps its own selection, just moved.
$ws1110.Range("B3").Select()

# "Thermo" becomes the active sheet/tab, selection on B6 (scrolled so row 6
# is at the top of the viewport).
$wsThermo.Activate()
$excel.ActiveWindow.ScrollRow = 6
$wsThermo.Range("B6").Select()

# "1220" keeps its previous selection (B4); it simply stops being the
# active/tabSelected sheet, which Activate() above already took care of.
